$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. ODI Batting sheet: rename MATCH_CARD_LINK -> MATCH_CODE header,
#    and replace the full scorecard URL in column D with just the
#    numeric MatchCode for every data row. Also drop the two stray
#    empty INNING_NUMBER cells (B21, B29).
# ------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"

for ($r = 2; $r -le 90; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $txt = $cell.Text
    if ($txt -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.Style = "Normal"
    }
}

$battingSheet.Range("B21").ClearContents()
$battingSheet.Range("B29").ClearContents()

# ------------------------------------------------------------------
# 2. ODI Bowling sheet: same MATCH_CARD_LINK -> MATCH_CODE rename,
#    this time it is column B.
# ------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"

for ($r = 2; $r -le 11; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $txt = $cell.Text
    if ($txt -match "MatchCode=(\d+)") {
        $code = $matches[1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.Style = "Normal"
    }
}

Write-Host "Stage 1+2 done"

# ------------------------------------------------------------------
# 3. New "Player Info" sheet, inserted before "ODI Batting".
# ------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le 4; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$idCell = $playerInfo.Cells.Item(2, 1)
$idCell.NumberFormat = "@"
$idCell.Value = "3607"
$idCell.Style = "Normal"

$playerInfo.Cells.Item(2, 2).Value = "Kyle James Coetzer"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

Write-Host "Stage 3 done"

# ------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet, appended after "ODI Bowling".
# ------------------------------------------------------------------
$bowlingSheetFresh = $wb.Worksheets.Item("ODI Bowling")
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $bowlingSheetFresh)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le 6; $c++) {
    $cell = $extra.Cells.Item(1, $c)
    $cell.Value = $exHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

Write-Host "Stage 4 header done"
